$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Update the diagram's wp14:editId token on the inline drawing in the
#    first paragraph. This value is not part of the document's visible
#    text (Find/Replace cannot reach it), so we rebuild that single
#    paragraph's XML with the new editId via InsertXML, leaving every
#    other attribute/relationship untouched.
# ---------------------------------------------------------------------
$drawingPara = $d.Paragraphs.Item(1)
$drawingFragment = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:noProof/><w:lang w:eastAsia="en-US"/></w:rPr><w:drawing><wp:inline xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" distT="0" distB="0" distL="0" distR="0" wp14:anchorId="4338F828" wp14:editId="6386DFA9"><wp:extent cx="5486400" cy="3200400"/><wp:effectExtent l="0" t="0" r="0" b="0"/><wp:docPr id="1" name="Diagram 1"/><wp:cNvGraphicFramePr/><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/diagram"><dgm:relIds xmlns:dgm="http://schemas.openxmlformats.org/drawingml/2006/diagram" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" r:dm="rId5" r:lo="rId6" r:qs="rId7" r:cs="rId8"/></a:graphicData></a:graphic></wp:inline></w:drawing></w:r></w:p>'
$drawingRange = $d.Range(0, $drawingPara.Range.End)
$drawingRange.InsertXML($drawingFragment)

# ---------------------------------------------------------------------
# 2) Remove the existing _GoBack bookmark. It will be re-created later,
#    inside the new final paragraph, right before the closing ")".
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------
# 3) Insert all of the new paragraphs (a blank separator followed by the
#    nine roles/responsibilities paragraphs) right at the document's
#    existing trailing empty paragraph. InsertXML replaces the target
#    range's content; by targeting a point just before that trailing
#    paragraph's own mark, every new paragraph we supply becomes a fresh
#    paragraph except the very last one, which merges into (and keeps
#    the identity of) that pre-existing trailing paragraph mark - i.e.
#    the trailing blank paragraph that used to sit right before the
#    sectPr effectively becomes the new "System engineering..." closing
#    paragraph, so no separate delete step is required.
# ---------------------------------------------------------------------
$wns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

$trailingPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$insertPoint = $d.Range($trailingPara.Range.End - 1, $trailingPara.Range.End - 1)

$newXml = "<w:p $wns/>" +
  "<w:p $wns>" +
    "<w:r><w:t>Roles</w:t></w:r>" +
    "<w:r><w:t xml:space='preserve'> and re</w:t></w:r>" +
    "<w:r><w:t xml:space='preserve'>sponsibilities </w:t></w:r>" +
    "<w:r><w:t>where more clarity of assignment is desired</w:t></w:r>" +
    "<w:r><w:t>:</w:t></w:r>" +
  "</w:p>" +
  "<w:p $wns>" +
    "<w:r><w:t>Satellite operations</w:t></w:r>" +
    "<w:r><w:t xml:space='preserve'> (unassigned)</w:t></w:r>" +
  "</w:p>" +
  "<w:p $wns>" +
    "<w:r><w:t>Customer support of Phase 4 radios</w:t></w:r>" +
    "<w:r><w:t xml:space='preserve'> (unassigned)</w:t></w:r>" +
  "</w:p>" +
  "<w:p $wns>" +
    "<w:r><w:t>Manufacturing of Phase 4 radios</w:t></w:r>" +
    "<w:r><w:t xml:space='preserve'> (unassigned)</w:t></w:r>" +
  "</w:p>" +
  "<w:p $wns>" +
    "<w:r><w:t>Development of 10GHz transverter</w:t></w:r>" +
    "<w:r><w:t xml:space='preserve'> (unassigned?)</w:t></w:r>" +
  "</w:p>" +
  "<w:p $wns>" +
    "<w:r><w:t>Development of firecode receiver</w:t></w:r>" +
    "<w:r><w:t xml:space='preserve'> (unassigned)</w:t></w:r>" +
  "</w:p>" +
  "<w:p $wns>" +
    "<w:r><w:t>Waveform development</w:t></w:r>" +
    "<w:r><w:t xml:space='preserve'> (currently Phase 4 Ground + VT</w:t></w:r>" +
    "<w:r><w:t>)</w:t></w:r>" +
  "</w:p>" +
  "<w:p $wns>" +
    "<w:r><w:t>Integration and Test of Ground vs. Satellite</w:t></w:r>" +
    "<w:r><w:t xml:space='preserve'> (unassigned, assumed Phase 4 Ground + VT</w:t></w:r>" +
    "<w:r><w:t>)</w:t></w:r>" +
  "</w:p>" +
  "<w:p $wns>" +
    "<w:r><w:t>System</w:t></w:r>" +
    "<w:r><w:t xml:space='preserve'> engineering (currently Phase 4 Ground + VT</w:t></w:r>" +
    "<w:r><w:t>)</w:t></w:r>" +
  "</w:p>"

$insertPoint.InsertXML($newXml)

# ---------------------------------------------------------------------
# 4) Re-add the _GoBack bookmark inside the new final paragraph, right
#    before its closing ")" run (i.e. right after "...Ground + VT").
# ---------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$bmPos = $lastPara.Range.End - 2
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)
